# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout/styling) and place the
#    copy right before "总计", renaming it "2022-Q1".
# 2. Fill the new sheet with the 2022-Q1 per-fund holdings.
# 3. Insert a new first data row in "总计" summarising the 2022-Q1 quarter,
#    pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1 & 2. New "2022-Q1" sheet (cloned from "2021-Q4" so header/styles match)
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$ws = $wb.Worksheets.Item(4)
$ws.Name = "2022-Q1"

# Columns B:G hold fund codes / numeric-looking values that must stay TEXT
# (inline/shared string) so leading zeros in fund codes (e.g. "013414") and
# the decimal-formatted figures survive, matching the source data. Force a
# text number format before writing, then strip the format back off so no
# stray "@" style lingers on the cells (values remain text either way).
$textRange = $ws.Range("B2:G5")
$textRange.NumberFormat = "@"

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "013414"
$ws.Cells.Item(2,3).Value = "太平智远三个月定期开放股票"
$ws.Cells.Item(2,4).Value = "8.69"
$ws.Cells.Item(2,5).Value = "86.34"
$ws.Cells.Item(2,6).Value = "8.70"
$ws.Cells.Item(2,7).Value = "0.7560"
$ws.Cells.Item(2,8).Value = 1

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "005270"
$ws.Cells.Item(3,3).Value = "太平改革红利精选灵活配置混合"
$ws.Cells.Item(3,4).Value = "1.87"
$ws.Cells.Item(3,5).Value = "88.32"
$ws.Cells.Item(3,6).Value = "8.70"
$ws.Cells.Item(3,7).Value = "0.1627"
$ws.Cells.Item(3,8).Value = 1

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "010896"
$ws.Cells.Item(4,3).Value = "太平价值增长股票A"
$ws.Cells.Item(4,4).Value = "1.18"
$ws.Cells.Item(4,5).Value = "83.63"
$ws.Cells.Item(4,6).Value = "8.85"
$ws.Cells.Item(4,7).Value = "0.1044"
$ws.Cells.Item(4,8).Value = 1

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "010897"
$ws.Cells.Item(5,3).Value = "太平价值增长股票C"
$ws.Cells.Item(5,4).Value = "1.01"
$ws.Cells.Item(5,5).Value = "83.63"
$ws.Cells.Item(5,6).Value = "8.85"
$ws.Cells.Item(5,7).Value = "0.0894"
$ws.Cells.Item(5,8).Value = 1

$textRange.ClearFormats()

# ---------------------------------------------------------------------------
# 3. "总计" sheet: push rows down one and insert the 2022-Q1 summary row
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Cells.Item(5,2).Value = "2020-Q4"
$totals.Cells.Item(5,3).Value = 4
$totals.Cells.Item(5,4).Value = 0.7

$totals.Cells.Item(4,2).Value = "2021-Q3"
$totals.Cells.Item(4,3).Value = 3
$totals.Cells.Item(4,4).Value = 0.36

$totals.Cells.Item(3,2).Value = "2021-Q4"
$totals.Cells.Item(3,3).Value = 4
$totals.Cells.Item(3,4).Value = 1.29

$totals.Cells.Item(2,2).Value = "2022-Q1"
$totals.Cells.Item(2,3).Value = 4
$totals.Cells.Item(2,4).Value = 1.11

# Column A is a running 0-based index with the same bold/bordered style as
# the header row; clone that style onto the newly-created A5 before writing
# its value (rows 2-4 already carry it from the original sheet).
$totals.Cells.Item(4,1).Copy()
$totals.Cells.Item(5,1).PasteSpecial(-4122)

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(5,1).Value = 3
